$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet 1: "Metadata"
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Metadata")

# First, create row 16 by copying the formatting (style) of row 15 so the
# new row reuses the existing body-row style (s="2") rather than creating
# a brand-new style entry in styles.xml.
$ws1.Range("A15:B15").Copy($ws1.Range("A16:B16"))

# Shift the content of rows 10-15 down to rows 11-16 (bottom-up so we
# never clobber a source row before it has been read). Using .Value
# (rather than .Copy) avoids a quirk where copying a truly blank cell
# does not clear the destination cell's existing content.
$ws1.Range("A16").Value = $ws1.Range("A15").Value
$ws1.Range("B16").Value = $ws1.Range("B15").Value

$ws1.Range("A15").Value = $ws1.Range("A14").Value
$ws1.Range("B15").Value = $ws1.Range("B14").Value

$ws1.Range("A14").Value = $ws1.Range("A13").Value
$ws1.Range("B14").Value = $ws1.Range("B13").Value

$ws1.Range("A13").Value = $ws1.Range("A12").Value
$ws1.Range("B13").Value = $ws1.Range("B12").Value

$ws1.Range("A12").Value = $ws1.Range("A11").Value
$ws1.Range("B12").Value = $ws1.Range("B11").Value

$ws1.Range("A11").Value = $ws1.Range("A10").Value
$ws1.Range("B11").Value = $ws1.Range("B10").Value

# Now overwrite the cells with their new final content.
$ws1.Range("B9").Value = "HL7 International / Cross-Group Projects"
$ws1.Range("B10").Value = "HL7 International / Cross-Group Projects (http://www.hl7.org/Special/committees/cgp, cgp@lists.HL7.org)"

$ws1.Range("A11").Value = "Contact"
$ws1.Range("B11").Value = "Health eData Inc (mailto:ehaas@healthedatainc.com)"

$ws1.Range("A12").Value = "Jurisdiction"
$ws1.Range("B12").Value = "United States of America"

$ws1.Range("A13").Value = "Description"
$ws1.Range("B13").Value = "Foo Bar Test ValueSet"

$ws1.Range("A14").Value = "Purpose"
$ws1.Range("B14").Value = ""

$ws1.Range("A15").Value = "Copyright"
$ws1.Range("B15").Value = "All rights reserved ... Don't mess with Greatness!!!"

$ws1.Range("A16").Value = "Immutable"
$ws1.Range("B16").Value = "BooleanType[null]"

# ----------------------------------------------------------------------
# Sheet 2: "Include from Foo Bar Test Cod"
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Include from Foo Bar Test Cod")

$ws2.Range("B1").Value = "Description"
